# Add two new data rows (53 and 54) to the "NEW" worksheet, matching the
# structure of the existing data rows (text columns A-L, O, P stored as
# plain text; numeric coordinate columns M, N stored as numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $value)
    # Prefix with an apostrophe so Excel treats the value as literal text
    # (prevents auto-conversion of numeric-looking / date-looking strings
    # into numbers or dates), then reset the style to "Normal" so the cell
    # does not pick up an implicit Text number format / style index.
    $ws.Range($cell).Value = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

function Set-NumberCell {
    param($cell, $value)
    $ws.Range($cell).Value = $value
}

# Row 53
Set-TextCell "A53" "6362"
Set-TextCell "B53" "7/8/2025"
Set-TextCell "C53" "ARIAS 3422"
Set-TextCell "D53" "12"
Set-TextCell "E53" "808099435"
Set-TextCell "F53" "NEW"
Set-TextCell "G53" "Pendiente"
Set-TextCell "H53" "Poste inclinado mal ubicado"
Set-TextCell "I53" "1"
Set-TextCell "J53" "Cambio"
Set-TextCell "K53" "Sin equipos"
Set-TextCell "L53" "Poste"
Set-NumberCell "M53" -58.483313
Set-NumberCell "N53" -34.54605
Set-TextCell "O53" "Saavedra"
Set-TextCell "P53" "Capital Norte"

# Row 54
Set-TextCell "A54" "6363"
Set-TextCell "B54" "7/8/2025"
Set-TextCell "C54" "MOLDES 3730"
Set-TextCell "D54" "12"
Set-TextCell "E54" "808099415"
Set-TextCell "F54" "NEW"
Set-TextCell "G54" "Pendiente"
Set-TextCell "H54" "Poste inclinado"
Set-TextCell "I54" "1"
Set-TextCell "J54" "Aplomo"
Set-TextCell "K54" "Sin equipos"
Set-TextCell "L54" "Poste"
Set-NumberCell "M54" -58.47192
Set-NumberCell "N54" -34.549398
Set-TextCell "O54" "Saavedra"
Set-TextCell "P54" "Capital Norte"

$ws.UsedRange.Address()
